$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-05-03T11:17:55+00:00"

# --- Elements sheet: update ID/Path/Base Path values + column widths ---
$elems = $wb.Worksheets.Item("Elements")

$elems.Range("A3").Value = "SourcePatientId.CX1"
$elems.Range("B3").Value = "SourcePatientId.CX1"
$elems.Range("AF3").Value = "SourcePatientId.CX1"

$elems.Range("A4").Value = "SourcePatientId.CX4"
$elems.Range("B4").Value = "SourcePatientId.CX4"
$elems.Range("AF4").Value = "SourcePatientId.CX4"

$elems.Range("A5").Value = "SourcePatientId.CX5"
$elems.Range("B5").Value = "SourcePatientId.CX5"
$elems.Range("AF5").Value = "SourcePatientId.CX5"

# The ID/Path/Base Path text grew by 2 chars ("CX" inserted), so Excel's
# best-fit autosize widens columns A, B and AF from 14.83203125 to
# 17.0078125 (stored width units). 16.1666... is the ColumnWidth (chars)
# value that maps to that stored width.
$elems.Columns.Item(1).ColumnWidth = 16.166666666666668
$elems.Columns.Item(2).ColumnWidth = 16.166666666666668
$elems.Columns.Item(32).ColumnWidth = 16.166666666666668
